$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()

$ws1.Range("A16").Value = "R"
$ws1.Range("B16").Value = 200000
$ws1.Range("C16").Value = "Ohm"
$ws1.Range("E16").Value = "Input impedance of the THS4531AIRUNR opamp"
$ws1.Range("E16").Font.Italic = $true

$ws1.Range("A17").Value = "C"
$ws1.Range("B17").Value = 0.0000001
$ws1.Range("B17").NumberFormat = "0.00E+00"
$ws1.Range("C17").Value = "F"
$ws1.Range("E17").Value = "Input DC blocking cap"
$ws1.Range("E17").Font.Italic = $true

$ws1.Range("A18").Value = "Fc"
$ws1.Range("B18").Formula = "=1/(2*PI()*B16*B17)"
$ws1.Range("C18").Value = "Hz"
$ws1.Range("E18").Font.Italic = $true

$ws1.Range("A15").Value = "Input differential 1st order HP filter"
$ws1.Range("A15").Font.Bold = $true

$ws1.Range("B18").Select() | Out-Null
